# Historial de pedidos de un cliente
# Remove the "Pelo de vaca" order row (row 10) from the Materias Primas sheet.
# Deleting the entire row shifts the rows below it up and Excel prunes the
# now-unused shared string ("Pelo de vaca") from the shared strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(10).Delete()
